# "Edited SPI part 2 ppt class 18"
#
# Slide 24 (sldId 328, cId 1648205850) has a body placeholder (shape 1,
# "Text Placeholder 1") with three paragraphs:
#   1) "SPI (Serial Peripheral Interface) is a synchronous serial ..."
#   2) (empty paragraph)
#   3) "It was developed by Motorola ..." (endParaRPr lang="en-IN")
#
# The author selected all of that body text and deleted it, leaving a
# single empty paragraph (keeping the lang="en-IN" of the last paragraph)
# with its bullet turned off.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$sp = $s.Shapes.Item(1)
$tr = $sp.TextFrame.TextRange

# Collapse the three paragraphs down to the last one: deleting paragraph 1
# merges it away (paragraph 2 shifts up to take its place), and repeating
# removes the former paragraph 2 as well, leaving only the original third
# paragraph (which keeps its own endParaRPr, lang="en-IN").
$tr.Paragraphs(1, 1).Delete()
$tr.Paragraphs(1, 1).Delete()

# Clear any remaining text in what is now the only paragraph.
$tr.Paragraphs(1, 1).Text = ""

# Turn the bullet off on the now-empty paragraph.
$tr.ParagraphFormat.Bullet.Visible = $false
